# "wednesday, finally the geomaps :)"
#
# The workbook's sheet "1" (Berlin districts table) is re-sorted so that its
# data rows (A3:F14) go back to ascending order by the district-code column
# (column A), instead of the previous descending-by-percentage (column E)
# sort. The last active cell/selection on that sheet also moved to F12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1")
$ws.Activate()

# Sort the data body (A3:F14) ascending by column A, leaving the header
# row (row 2) and the autofilter untouched.
$sortRange = $ws.Range("A3:F14")
$keyColumn = $ws.Range("A3:A14")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyColumn, 0, 1, 0, 0) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

# Reflect the author's final selection on the sheet.
$ws.Range("F12").Select() | Out-Null
